# The document's title/description were renamed from "Personal Open Science
# Pipeline" / "POSP" to "Open Science Pipeline" / "OSP". This touches the
# heading text, the intro paragraph text, and the heading's bookmark name.
# (The Bookmarks collection's Name-assignment / Delete+Add round trip is not
# reliable in this host for in-place renames without disturbing paragraph
# structure, so the safest, most faithful way to apply these small textual
# edits - including the bookmark *name* attribute, which isn't reachable via
# Find/Replace since it isn't part of the document's visible text - is a
# direct round trip through the document's WordOpenXML package.)

$d = $word.ActiveDocument

$xml = $d.WordOpenXML

$xml = $xml.Replace('w:name="personal-open-science-pipeline"', 'w:name="open-science-pipeline"')
$xml = $xml.Replace('Personal Open Science Pipeline', 'Open Science Pipeline')
$xml = $xml.Replace('POSP aims to develop', 'OSP aims to develop')
$xml = $xml.Replace('6f191523', 'eddda7b7')

$d.Content.InsertXML($xml) | Out-Null
